$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value2
$text = $text -replace [regex]::Escape("1000 Bs = 6.68 = 26669.25 pesos"), "1000 Bs = 6.63 = 26458.89 pesos"
$text = $text -replace [regex]::Escape("26669.25 pesos = 6.65 = 957.7 Bs"), "26458.89 pesos = 6.62 = 973.47 Bs"
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 150.8
$wsTasas.Range("O10").Value = 3990
$wsTasas.Range("N12").Value = 3995.99
$wsTasas.Range("O12").Value = 147.02
